# Adds four new worksheets (test10..test13), each a copy of the existing
# "test9" sheet, to exercise "find the used cell range of a worksheet"
# against a variety of selections:
#   test10 -> edits C4, then selects single cell C4
#   test11 -> edits C4, then selects a single cell far outside the data (E24)
#   test12 -> edits E8,  then selects a multi-cell range outside the data (L1:P4)
#   test13 -> edits E8,  then selects the last real data cell (E8) and ends up active

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("test9")

# --- test10 ---------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws10 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws10.Name = "test10"
$ws10.Range("C4").Value = "blah10"
$ws10.Range("C4").Select()

# --- test11 ---------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws11 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws11.Name = "test11"
$ws11.Range("C4").Value = "blah11"
$ws11.Range("E24").Select()

# --- test12 ---------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws12 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws12.Name = "test12"
$ws12.Range("E8").Value = "blah12"
$ws12.Range("L1:P4").Select()

# --- test13 ---------------------------------------------------------------
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws13 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws13.Name = "test13"
$ws13.Range("E8").Value = "blah13"
$ws13.Range("E8").Select()

# test13 is now the active sheet/tab (last one activated via the rename/select
# above), matching the workbook's new activeTab.
